$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cells for existing cell styles already present in the workbook:
#   A2 -> default / no fill style
#   D2 -> green fill style (used for "SI"/passing values)
#   H2 -> red fill style (used for "NO"/failing values)
$noStyleRef    = $ws.Range("A2")
$greenStyleRef = $ws.Range("D2")
$redStyleRef   = $ws.Range("H2")

function Set-GreenNumber {
    param($cellRef, $value)
    $greenStyleRef.Copy()
    $ws.Range($cellRef).PasteSpecial(-4122)  # xlPasteFormats
    $ws.Range($cellRef).Value = $value
}

function Set-NoStyleText {
    param($cellRef, $text)
    $noStyleRef.Copy()
    $ws.Range($cellRef).PasteSpecial(-4122)  # xlPasteFormats
    $ws.Range($cellRef).Value = $text
}

function Set-RedText {
    param($cellRef, $text)
    $redStyleRef.Copy()
    $ws.Range($cellRef).PasteSpecial(-4122)  # xlPasteFormats
    $ws.Range($cellRef).Value = $text
}

# ---- Rows 2-4: L/M/N/O go from "N/A" to numeric scores (green fill) ----
foreach ($r in 2,3,4) {
    Set-GreenNumber "L$r" 96
    Set-GreenNumber "M$r" 99
    Set-GreenNumber "N$r" 100
    Set-GreenNumber "O$r" 100
}

# ---- Rows 5-7: L/M go from "N/A" to numeric scores (N and O remain "N/A") ----
foreach ($r in 5,6,7) {
    Set-GreenNumber "L$r" 91
    Set-GreenNumber "M$r" 91
}

# ---- Rows 11-13: L/M/N/O go from "N/A" to numeric scores ----
foreach ($r in 11,12,13) {
    Set-GreenNumber "L$r" 91
    Set-GreenNumber "M$r" 97
    Set-GreenNumber "N$r" 90
    Set-GreenNumber "O$r" 100
}

# ---- Rows 14-16: L/M/N/O go from "N/A" to numeric scores ----
foreach ($r in 14,15,16) {
    Set-GreenNumber "L$r" 91
    Set-GreenNumber "M$r" 98
    Set-GreenNumber "N$r" 90
    Set-GreenNumber "O$r" 100
}

# ---- Row 17: L/M/N/O go from "N/A" to numeric scores ----
Set-GreenNumber "L17" 91
Set-GreenNumber "M17" 99
Set-GreenNumber "N17" 90
Set-GreenNumber "O17" 100

# ---- Row 18: page now errors out; text fields become ERROR, checks flip to NO,
#      "0 de 10" -> "0 de 0" with default style, L-O become numeric scores ----
Set-NoStyleText "C18" "ERROR"
Set-RedText     "D18" "NO"
Set-NoStyleText "E18" "ERROR"
Set-RedText     "F18" "NO"
Set-NoStyleText "G18" "ERROR"
Set-NoStyleText "I18" "ERROR"
Set-NoStyleText "J18" "0 de 0"
Set-GreenNumber "L18" 91
Set-GreenNumber "M18" 99
Set-GreenNumber "N18" 90
Set-GreenNumber "O18" 100

# ---- Row 19 ----
Set-NoStyleText "C19" "ERROR"
Set-NoStyleText "E19" "ERROR"
Set-NoStyleText "G19" "ERROR"
Set-NoStyleText "I19" "ERROR"
Set-NoStyleText "J19" "0 de 0"
Set-GreenNumber "L19" 91
Set-GreenNumber "M19" 99
Set-GreenNumber "N19" 90
Set-GreenNumber "O19" 100

# ---- Row 20 (L-O stay "N/A", not touched) ----
Set-NoStyleText "C20" "ERROR"
Set-RedText     "D20" "NO"
Set-NoStyleText "E20" "ERROR"
Set-NoStyleText "G20" "ERROR"
Set-NoStyleText "I20" "ERROR"
Set-NoStyleText "J20" "0 de 0"

# ---- Row 21 (L-O stay "N/A", not touched) ----
Set-NoStyleText "C21" "ERROR"
Set-RedText     "D21" "NO"
Set-NoStyleText "E21" "ERROR"
Set-RedText     "F21" "NO"
Set-NoStyleText "G21" "ERROR"
Set-NoStyleText "I21" "ERROR"
Set-NoStyleText "J21" "0 de 0"

# ---- Row 22 (L-O stay "N/A", not touched) ----
Set-NoStyleText "C22" "ERROR"
Set-NoStyleText "E22" "ERROR"
Set-NoStyleText "G22" "ERROR"
Set-NoStyleText "I22" "ERROR"
Set-NoStyleText "J22" "0 de 0"

$wb.Application.CutCopyMode = $false
